$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Estimated")
$ws.Range("B2").Value = "Mon Oct 02 17:51:41 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 17:52:22 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 17:53:01 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 02 17:53:38 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 02 17:54:16 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 02 17:54:55 EDT 2023"

$ws = $wb.Worksheets.Item("Existing")
$ws.Range("B2").Value = "Mon Oct 02 17:55:33 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 17:56:11 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 17:56:50 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 02 17:57:28 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 02 17:58:06 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 02 17:58:44 EDT 2023"
$ws.Range("B8").Value = "Mon Oct 02 17:59:23 EDT 2023"
$ws.Range("B9").Value = "Mon Oct 02 18:00:01 EDT 2023"
$ws.Range("B10").Value = "Mon Oct 02 18:00:41 EDT 2023"
$ws.Range("B11").Value = "Mon Oct 02 18:01:19 EDT 2023"
$ws.Range("B12").Value = "Mon Oct 02 18:01:57 EDT 2023"

$ws = $wb.Worksheets.Item("Extension")
$ws.Range("B2").Value = "Mon Oct 02 18:02:36 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 18:03:15 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 18:03:51 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 02 18:04:27 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 02 18:05:03 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 02 18:05:38 EDT 2023"

$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Range("B2").Value = "Mon Oct 02 18:06:15 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 18:06:51 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 18:07:27 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 02 18:08:03 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 02 18:08:40 EDT 2023"
$ws.Range("B7").Value = "Mon Oct 02 18:09:15 EDT 2023"
$ws.Range("B8").Value = "Mon Oct 02 18:09:53 EDT 2023"
$ws.Range("B9").Value = "Mon Oct 02 18:10:30 EDT 2023"
$ws.Range("B10").Value = "Mon Oct 02 18:11:07 EDT 2023"
$ws.Range("B11").Value = "Mon Oct 02 18:11:44 EDT 2023"
$ws.Range("B12").Value = "Mon Oct 02 18:12:21 EDT 2023"
$ws.Range("B13").Value = "Mon Oct 02 18:12:57 EDT 2023"
$ws.Range("B14").Value = "Mon Oct 02 18:13:33 EDT 2023"
$ws.Range("B15").Value = "Mon Oct 02 18:14:10 EDT 2023"
$ws.Range("B16").Value = "Mon Oct 02 18:14:47 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Range("B2").Value = "Mon Oct 02 18:15:59 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 18:16:32 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 18:17:07 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 02 18:17:41 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 02 18:18:15 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Range("B2").Value = "Mon Oct 02 18:18:50 EDT 2023"
$ws.Range("B3").Value = "Mon Oct 02 18:19:31 EDT 2023"
$ws.Range("B4").Value = "Mon Oct 02 18:20:12 EDT 2023"
$ws.Range("B5").Value = "Mon Oct 02 18:20:54 EDT 2023"
$ws.Range("B6").Value = "Mon Oct 02 18:21:35 EDT 2023"

$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Range("B2").Value = "Mon Oct 02 18:15:23 EDT 2023"
